$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.525.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.01%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.919.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.39%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.81%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'325.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.14%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.71%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4813"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.76%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.4057"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.48%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.08222"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.98%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'1.010"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.31%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'23.40"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.38%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.940.70"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.00%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'6.056"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.39%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'7.239"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.85%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'91.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.26%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.06862"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.03%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.014"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.64%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.00001039"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.21%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'17.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.91%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'1.012"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.71%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'29.527.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.04%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.673"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.10%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.46%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +1.18%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.152.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.97%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'6.537"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +4.11%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'155.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.68%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'20.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.05%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.099"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.33%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'120.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.58%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.018"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.56%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.09633"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.69%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'5.619"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.43%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.17%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.373"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.68%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.06313"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +3.20%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.02283"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.65%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.182"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.79%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.5937"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.22%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'10.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.11%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'7.903"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.23%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1848"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.40%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.465"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.79%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.282"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.37%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'12.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.23%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.07471"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.28%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.5552"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.42%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.942"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.67%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'118.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.06%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.431"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.35%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'72.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.82%  "
$ws.Range("E51").Style = "Normal"
